$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83
$ws.Cells.Item(83,1).Value = 82
$ws.Cells.Item(83,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(83,3).Value = "8:50 PM"
$ws.Cells.Item(83,4).Value = "LO3825"
$ws.Cells.Item(83,5).Value = "Warsaw"
$ws.Cells.Item(83,6).Value = "(WAW)"
$ws.Cells.Item(83,7).Value = "LOT "
$ws.Cells.Item(83,8).Value = "E75S"
$ws.Cells.Item(83,9).Value = "(SP-LIQ)"
$ws.Cells.Item(83,10).Value = "8:46 PM"
$ws.Cells.Item(83,12).Value = "0 hours, -4 minutes"

# Row 84
$ws.Cells.Item(84,1).Value = 83
$ws.Cells.Item(84,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(84,3).Value = "8:55 PM"
$ws.Cells.Item(84,4).Value = "LH1646"
$ws.Cells.Item(84,5).Value = "Munich"
$ws.Cells.Item(84,6).Value = "(MUC)"
$ws.Cells.Item(84,7).Value = "Lufthansa "
$ws.Cells.Item(84,8).Value = "CRJ9"
$ws.Cells.Item(84,9).Value = "(D-ACNW)"
$ws.Cells.Item(84,10).Value = "8:52 PM"
$ws.Cells.Item(84,12).Value = "0 hours, -3 minutes"

# Row 85
$ws.Cells.Item(85,1).Value = 84
$ws.Cells.Item(85,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(85,3).Value = "9:15 PM"
$ws.Cells.Item(85,4).Value = "FR4616"
$ws.Cells.Item(85,5).Value = "Stockholm"
$ws.Cells.Item(85,6).Value = "(ARN)"
$ws.Cells.Item(85,7).Value = "Ryanair "
$ws.Cells.Item(85,8).Value = "B738"
$ws.Cells.Item(85,9).Value = "(9H-QBG)"
$ws.Cells.Item(85,10).Value = "9:44 PM"
$ws.Cells.Item(85,12).Value = "0 hours, 29 minutes"

# Row 86
$ws.Cells.Item(86,1).Value = 85
$ws.Cells.Item(86,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(86,3).Value = "9:55 PM"
$ws.Cells.Item(86,4).Value = "W61762"
$ws.Cells.Item(86,5).Value = "Trondheim"
$ws.Cells.Item(86,6).Value = "(TRD)"
$ws.Cells.Item(86,7).Value = "Wizz Air "
$ws.Cells.Item(86,8).Value = "A21N"
$ws.Cells.Item(86,9).Value = "(9H-WAP)"
$ws.Cells.Item(86,10).Value = "9:47 PM"
$ws.Cells.Item(86,12).Value = "0 hours, -8 minutes"

# Row 87
$ws.Cells.Item(87,1).Value = 86
$ws.Cells.Item(87,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(87,3).Value = "10:10 PM"
$ws.Cells.Item(87,4).Value = "W61676"
$ws.Cells.Item(87,5).Value = "Hamburg"
$ws.Cells.Item(87,6).Value = "(HAM)"
$ws.Cells.Item(87,7).Value = "Wizz Air "
$ws.Cells.Item(87,8).Value = "A320"
$ws.Cells.Item(87,9).Value = "(HA-LWV)"
$ws.Cells.Item(87,10).Value = "9:41 PM"
$ws.Cells.Item(87,12).Value = "0 hours, -29 minutes"

# Row 88
$ws.Cells.Item(88,1).Value = 87
$ws.Cells.Item(88,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(88,3).Value = "10:15 PM"
$ws.Cells.Item(88,4).Value = "FR6120"
$ws.Cells.Item(88,5).Value = "Hamburg"
$ws.Cells.Item(88,6).Value = "(HAM)"
$ws.Cells.Item(88,7).Value = "Ryanair "
$ws.Cells.Item(88,8).Value = "B738"
$ws.Cells.Item(88,9).Value = "(SP-RKQ)"
$ws.Cells.Item(88,10).Value = "10:05 PM"
$ws.Cells.Item(88,12).Value = "0 hours, -10 minutes"

# Row 89
$ws.Cells.Item(89,1).Value = 88
$ws.Cells.Item(89,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(89,3).Value = "10:35 PM"
$ws.Cells.Item(89,4).Value = "FR2592"
$ws.Cells.Item(89,5).Value = "Valencia"
$ws.Cells.Item(89,6).Value = "(VLC)"
$ws.Cells.Item(89,7).Value = "Ryanair "
$ws.Cells.Item(89,8).Value = "B738"
$ws.Cells.Item(89,9).Value = "(SP-RSW)"
$ws.Cells.Item(89,10).Value = "10:34 PM"
$ws.Cells.Item(89,12).Value = "0 hours, -1 minutes"

# Row 90
$ws.Cells.Item(90,1).Value = 89
$ws.Cells.Item(90,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(90,3).Value = "10:55 PM"
$ws.Cells.Item(90,4).Value = "W61734"
$ws.Cells.Item(90,5).Value = "Stockholm"
$ws.Cells.Item(90,6).Value = "(NYO)"
$ws.Cells.Item(90,7).Value = "Wizz Air "
$ws.Cells.Item(90,8).Value = "A321"
$ws.Cells.Item(90,9).Value = "(HA-LTB)"
$ws.Cells.Item(90,10).Value = "10:37 PM"
$ws.Cells.Item(90,12).Value = "0 hours, -18 minutes"

Write-Host ("Dimension: " + $ws.UsedRange.Address())

